# cambio para evaluar en una sola fecha y que se pueda actualizar todo el dataset de una
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded second date row ("fecha2") entirely -- everything
# below shifts up one row.
$ws.Rows.Item(3).EntireRow.Delete()

# Row 2 ("fecha1" -> "fecha"): single evaluation date, updated value.
$ws.Range("A2").Value = "fecha"
$ws.Range("C2").Value = "29_01_2025"

# Row 3 (was row 4): DNI updated to the new dataset's identifier.
$ws.Range("C3").Value = 34483112

# Row 6 (was row 7): new objective text for the dataset refresh.
$ws.Range("C6").Value = "Escalar 8b+ en marzo comodo"

# Park the selection on the DNI cell, as in the refreshed template.
$ws.Range("C4").Select()
